# Generate Report for Archive
#
# The localization report is regenerated: every "Ready for handoff" status
# cell becomes "In Translation". Because the new text is shorter, the
# status column(s) that had been sized to fit the old text are narrowed
# to fit the new, shorter text (column width here mirrors what
# Columns.AutoFit() would produce in Excel for the new content).

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- 1. Replace the status text everywhere it appears -----------------
# (NOTE: keep $oldStatus on the left of -eq. Some cells hold a Boolean
#  True/False value, and PowerShell's -eq coerces the right-hand side to
#  match the left operand's type; "$cell.Value() -eq $oldStatus" would
#  therefore wrongly report a match for every True cell.)
foreach ($ws in $wb.Worksheets()) {
    $used = $ws.UsedRange()
    foreach ($cell in $used.Cells()) {
        if ($oldStatus -eq $cell.Value()) {
            $cell.Value = $newStatus
        }
    }
}

# --- 2. Re-fit the status columns to the new (shorter) text -----------
$newStatusColumnWidth = 12.5

# Overview sheet: status is duplicated in columns E (zh-cn) and F (de-de)
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).AutoFit()
$overview.Columns.Item(5).ColumnWidth = $newStatusColumnWidth
$overview.Columns.Item(6).AutoFit()
$overview.Columns.Item(6).ColumnWidth = $newStatusColumnWidth

# Per-language sheets: status lives in column C
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).AutoFit()
$zhcn.Columns.Item(3).ColumnWidth = $newStatusColumnWidth

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).AutoFit()
$dede.Columns.Item(3).ColumnWidth = $newStatusColumnWidth
